$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate the "UC_Sets" input block (rows 12:18, cols B:J) down to
#    rows 21:27, carrying over both values/formulas and formatting.
# ---------------------------------------------------------------------------
$ws.Range("B12:J18").Copy($ws.Range("B21"))

# ---------------------------------------------------------------------------
# 2. Fix up the copied block so it describes the new CO2EQS growth-uc
#    constraint instead of being a literal duplicate of the CO2S one.
# ---------------------------------------------------------------------------

# Row 21/22 (the "~UC_Sets" header lines) and row 23 (the "~UC_T" line) only
# had a couple of populated cells in the source rows (12/13/14); the
# rectangular Copy() above also materialised the in-between blanks as empty
# cells, so drop those back out to match the sparser source rows.
$ws.Range("C21").ClearContents()
$ws.Range("E21:J21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("F22:J22").ClearContents()
$ws.Range("B23:F23").ClearContents()
$ws.Range("H23:J23").ClearContents()

# Row 24 (the column-header line) is identical to row 15, so nothing else to
# do there - the copy already produced the right content.

# Row 25 (was row 16): new UC name + literal commodity string instead of a
# formula, and the "-E.." reference shifts down with the block.
$ws.Range("B25").Value = "UC_Growth_CO2EQS"
$ws.Range("D25").Value = "CO2EQS"
$ws.Range("I25").Formula = '=-E21'

# Row 26 (was row 17)
$ws.Range("B26").Value = "UC_Growth_CO2EQS"
$ws.Range("D26").Value = "CO2EQS"
$ws.Range("E26").Value = 2029
$ws.Range("I26").Formula = '=-E22'
# Entering a formula that references E22 (styled) leaks that cell's number
# format onto I26; put I26 back to the unstyled "General" look the diff
# expects (matches I16/I17/I18, none of which carry an "s" attribute).
$ws.Range("I26").Style = "Normal"

# Row 27 (was row 18)
$ws.Range("B27").Value = "UC_Growth_CO2EQS"
$ws.Range("D27").Value = "CO2EQS"
$ws.Range("E27").Value = 2030
$ws.Range("G27").Value = 1
$ws.Range("I27").Formula = '=-E22'
$ws.Range("I27").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Move the old chart-source data (growth rate + year headers + CO2 curve)
#    six columns to the right (I:N -> O:T) to make room, then rebuild the
#    formulas for the curve row without the old "shared formula" grouping.
# ---------------------------------------------------------------------------

# Growth-rate factor used by the CO2 decay curve: I25 -> O25
# (I25 itself was already overwritten above with the "=-E21" formula that
# belongs to the duplicated input block; J25 already holds the correct
# copied "5" from the block copy, so only K25:N25 need clearing.)
$ws.Range("O25").Value = 0.75
$ws.Range("K25:N25").ClearContents()

# Year headers: I27:N27 -> O27:T27
# (I27 already holds the "=-E22" formula and J27 already holds the copied
# "5" from the block copy, so only K27:N27 need clearing.)
$ws.Range("O27").Value = 2050
$ws.Range("P27").Value = 2051
$ws.Range("Q27").Value = 2052
$ws.Range("R27").Value = 2053
$ws.Range("S27").Value = 2054
$ws.Range("T27").Value = 2055
$ws.Range("K27:N27").ClearContents()

# CO2 curve: H28:N28 -> N28:T28 (clear the old block first, it overlaps the
# new N28 target cell)
$ws.Range("O28").Value = 220
$ws.Range("P28").Formula = '=O28*$O$25'
$ws.Range("Q28").Formula = '=P28*$O$25'
$ws.Range("R28").Formula = '=Q28*$O$25'
$ws.Range("S28").Formula = '=R28*$O$25'
$ws.Range("T28").Value = 0
$ws.Range("H28:N28").ClearContents()
$ws.Range("N28").Value = "CO2"

# ---------------------------------------------------------------------------
# 4. Point the chart at the relocated ranges.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = '=SERIES(Sheet1!$N$28,Sheet1!$O$27:$T$27,Sheet1!$O$28:$T$28,1)'

# Reflect where the author ended up looking after adding the new block.
[void]$ws.Range("C31").Select()

$wb.Application.Calculate()
